$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.1952936666666667
$ws.Cells.Item(2, 8).Value = 0.585881
$ws.Cells.Item(2, 9).Value = 0.001827617096392301
$ws.Cells.Item(2, 10).Value = 0.0018276170963923
$ws.Cells.Item(2, 13).Value = 42.04602466666667
$ws.Cells.Item(2, 14).Value = 126.138074
$ws.Cells.Item(2, 15).Value = 0.1180439555498783
$ws.Cells.Item(2, 16).Value = 0.1180439555498783
$ws.Cells.Item(2, 17).Value = 8.211322325910444
$ws.Cells.Item(2, 18).Value = 73.901900933194
$ws.Cells.Item(2, 19).Value = 0.0002157391512887303
$ws.Cells.Item(2, 20).Value = 0.0002157391512887302

$ws.Cells.Item(3, 7).Value = 0.1952936666666667
$ws.Cells.Item(3, 8).Value = 0.585881
$ws.Cells.Item(3, 9).Value = 0.001827617096392301
$ws.Cells.Item(3, 10).Value = 0.0018276170963923
$ws.Cells.Item(3, 15).Value = 0.1482760805823429
$ws.Cells.Item(3, 16).Value = 0.1482760805823429
$ws.Cells.Item(3, 17).Value = 10.31431626645067
$ws.Cells.Item(3, 18).Value = 92.828846398056
$ws.Cells.Item(3, 19).Value = 0.0002709918998583324
$ws.Cells.Item(3, 20).Value = 0.0002709918998583323

$ws.Cells.Item(4, 7).Value = 0.1952936666666667
$ws.Cells.Item(4, 8).Value = 0.585881
$ws.Cells.Item(4, 9).Value = 0.001827617096392301
$ws.Cells.Item(4, 10).Value = 0.0018276170963923
$ws.Cells.Item(4, 13).Value = 68.81807333333334
$ws.Cells.Item(4, 14).Value = 206.45422
$ws.Cells.Item(4, 15).Value = 0.1932063174578422
$ws.Cells.Item(4, 16).Value = 0.1932063174578422
$ws.Cells.Item(4, 17).Value = 13.43973387420223
$ws.Cells.Item(4, 18).Value = 120.95760486782
$ws.Cells.Item(4, 19).Value = 0.0003531071689169506
$ws.Cells.Item(4, 20).Value = 0.0003531071689169505

$ws.Cells.Item(5, 7).Value = 0.1952936666666667
$ws.Cells.Item(5, 8).Value = 0.585881
$ws.Cells.Item(5, 9).Value = 0.001827617096392301
$ws.Cells.Item(5, 10).Value = 0.0018276170963923
$ws.Cells.Item(5, 13).Value = 11.78107633333333
$ws.Cells.Item(5, 14).Value = 35.343229
$ws.Cells.Item(5, 15).Value = 0.03307529931894448
$ws.Cells.Item(5, 16).Value = 0.03307529931894448
$ws.Cells.Item(5, 17).Value = 2.300769594416555
$ws.Cells.Item(5, 18).Value = 20.706926349749
$ws.Cells.Item(5, 19).Value = 0.00006044898250359554
$ws.Cells.Item(5, 20).Value = 0.00006044898250359552

$ws.Cells.Item(6, 7).Value = 0.1952936666666667
$ws.Cells.Item(6, 8).Value = 0.585881
$ws.Cells.Item(6, 9).Value = 0.001827617096392301
$ws.Cells.Item(6, 10).Value = 0.0018276170963923
$ws.Cells.Item(6, 13).Value = 139.820737
$ws.Cells.Item(6, 14).Value = 419.462211
$ws.Cells.Item(6, 15).Value = 0.3925458588351179
$ws.Cells.Item(6, 16).Value = 0.3925458588351179
$ws.Cells.Item(6, 17).Value = 27.30610440476567
$ws.Cells.Item(6, 18).Value = 245.754939642891
$ws.Cells.Item(6, 19).Value = 0.0007174235227250601
$ws.Cells.Item(6, 20).Value = 0.0007174235227250599

$ws.Cells.Item(7, 7).Value = 0.1952936666666667
$ws.Cells.Item(7, 8).Value = 0.585881
$ws.Cells.Item(7, 9).Value = 0.001827617096392301
$ws.Cells.Item(7, 10).Value = 0.0018276170963923
$ws.Cells.Item(7, 13).Value = 40.90925733333334
$ws.Cells.Item(7, 14).Value = 122.727772
$ws.Cells.Item(7, 15).Value = 0.1148524882558742
$ws.Cells.Item(7, 16).Value = 0.1148524882558742
$ws.Cells.Item(7, 17).Value = 7.98931886523689
$ws.Cells.Item(7, 18).Value = 71.903869787132
$ws.Cells.Item(7, 19).Value = 0.0002099063710996316
$ws.Cells.Item(7, 20).Value = 0.0002099063710996316

$ws.Cells.Item(8, 7).Value = 3.363724333333333
$ws.Cells.Item(8, 8).Value = 10.091173
$ws.Cells.Item(8, 9).Value = 0.03147874789838274
$ws.Cells.Item(8, 10).Value = 0.03147874789838274
$ws.Cells.Item(8, 13).Value = 42.04602466666667
$ws.Cells.Item(8, 14).Value = 126.138074
$ws.Cells.Item(8, 15).Value = 0.1180439555498783
$ws.Cells.Item(8, 16).Value = 0.1180439555498783
$ws.Cells.Item(8, 17).Value = 141.4312362912002
$ws.Cells.Item(8, 18).Value = 1272.881126620802
$ws.Cells.Item(8, 19).Value = 0.003715875917682516
$ws.Cells.Item(8, 20).Value = 0.003715875917682515

$ws.Cells.Item(9, 7).Value = 3.363724333333333
$ws.Cells.Item(9, 8).Value = 10.091173
$ws.Cells.Item(9, 9).Value = 0.03147874789838274
$ws.Cells.Item(9, 10).Value = 0.03147874789838274
$ws.Cells.Item(9, 15).Value = 0.1482760805823429
$ws.Cells.Item(9, 16).Value = 0.1482760805823429
$ws.Cells.Item(9, 17).Value = 177.6530555206053
$ws.Cells.Item(9, 18).Value = 1598.877499685448
$ws.Cells.Item(9, 19).Value = 0.004667545360011858
$ws.Cells.Item(9, 20).Value = 0.004667545360011857

$ws.Cells.Item(10, 7).Value = 3.363724333333333
$ws.Cells.Item(10, 8).Value = 10.091173
$ws.Cells.Item(10, 9).Value = 0.03147874789838274
$ws.Cells.Item(10, 10).Value = 0.03147874789838274
$ws.Cells.Item(10, 13).Value = 68.81807333333334
$ws.Cells.Item(10, 14).Value = 206.45422
$ws.Cells.Item(10, 15).Value = 0.1932063174578422
$ws.Cells.Item(10, 16).Value = 0.1932063174578422
$ws.Cells.Item(10, 17).Value = 231.4850278444511
$ws.Cells.Item(10, 18).Value = 2083.36525060006
$ws.Cells.Item(10, 19).Value = 0.00608189295963032
$ws.Cells.Item(10, 20).Value = 0.006081892959630317

$ws.Cells.Item(11, 7).Value = 3.363724333333333
$ws.Cells.Item(11, 8).Value = 10.091173
$ws.Cells.Item(11, 9).Value = 0.03147874789838274
$ws.Cells.Item(11, 10).Value = 0.03147874789838274
$ws.Cells.Item(11, 13).Value = 11.78107633333333
$ws.Cells.Item(11, 14).Value = 35.343229
$ws.Cells.Item(11, 15).Value = 0.03307529931894448
$ws.Cells.Item(11, 16).Value = 0.03307529931894448
$ws.Cells.Item(11, 17).Value = 39.62829313529078
$ws.Cells.Item(11, 18).Value = 356.654638217617
$ws.Cells.Item(11, 19).Value = 0.001041169008924604
$ws.Cells.Item(11, 20).Value = 0.001041169008924603

$ws.Cells.Item(12, 7).Value = 3.363724333333333
$ws.Cells.Item(12, 8).Value = 10.091173
$ws.Cells.Item(12, 9).Value = 0.03147874789838274
$ws.Cells.Item(12, 10).Value = 0.03147874789838274
$ws.Cells.Item(12, 13).Value = 139.820737
$ws.Cells.Item(12, 14).Value = 419.462211
$ws.Cells.Item(12, 15).Value = 0.3925458588351179
$ws.Cells.Item(12, 16).Value = 0.3925458588351179
$ws.Cells.Item(12, 17).Value = 470.3184153515003
$ws.Cells.Item(12, 18).Value = 4232.865738163503
$ws.Cells.Item(12, 19).Value = 0.01235685212882482
$ws.Cells.Item(12, 20).Value = 0.01235685212882481

$ws.Cells.Item(13, 7).Value = 3.363724333333333
$ws.Cells.Item(13, 8).Value = 10.091173
$ws.Cells.Item(13, 9).Value = 0.03147874789838274
$ws.Cells.Item(13, 10).Value = 0.03147874789838274
$ws.Cells.Item(13, 13).Value = 40.90925733333334
$ws.Cells.Item(13, 14).Value = 122.727772
$ws.Cells.Item(13, 15).Value = 0.1148524882558742
$ws.Cells.Item(13, 16).Value = 0.1148524882558742
$ws.Cells.Item(13, 17).Value = 137.6074643507285
$ws.Cells.Item(13, 18).Value = 1238.467179156556
$ws.Cells.Item(13, 19).Value = 0.00361541252330863
$ws.Cells.Item(13, 20).Value = 0.003615412523308628

$ws.Cells.Item(14, 7).Value = 101.145495
$ws.Cells.Item(14, 8).Value = 303.436485
$ws.Cells.Item(14, 9).Value = 0.9465500804006033
$ws.Cells.Item(14, 10).Value = 0.9465500804006032
$ws.Cells.Item(14, 13).Value = 42.04602466666667
$ws.Cells.Item(14, 14).Value = 126.138074
$ws.Cells.Item(14, 15).Value = 0.1180439555498783
$ws.Cells.Item(14, 16).Value = 0.1180439555498783
$ws.Cells.Item(14, 17).Value = 4252.76597769221
$ws.Cells.Item(14, 18).Value = 38274.89379922989
$ws.Cells.Item(14, 19).Value = 0.1117345156165425
$ws.Cells.Item(14, 20).Value = 0.1117345156165425

$ws.Cells.Item(15, 7).Value = 101.145495
$ws.Cells.Item(15, 8).Value = 303.436485
$ws.Cells.Item(15, 9).Value = 0.9465500804006033
$ws.Cells.Item(15, 10).Value = 0.9465500804006032
$ws.Cells.Item(15, 15).Value = 0.1482760805823429
$ws.Cells.Item(15, 16).Value = 0.1482760805823429
$ws.Cells.Item(15, 17).Value = 5341.937821964039
$ws.Cells.Item(15, 18).Value = 48077.44039767636
$ws.Cells.Item(15, 19).Value = 0.140350735996703
$ws.Cells.Item(15, 20).Value = 0.140350735996703

$ws.Cells.Item(16, 7).Value = 101.145495
$ws.Cells.Item(16, 8).Value = 303.436485
$ws.Cells.Item(16, 9).Value = 0.9465500804006033
$ws.Cells.Item(16, 10).Value = 0.9465500804006032
$ws.Cells.Item(16, 13).Value = 68.81807333333334
$ws.Cells.Item(16, 14).Value = 206.45422
$ws.Cells.Item(16, 15).Value = 0.1932063174578422
$ws.Cells.Item(16, 16).Value = 0.1932063174578422
$ws.Cells.Item(16, 17).Value = 6960.638092246301
$ws.Cells.Item(16, 18).Value = 62645.74283021671
$ws.Cells.Item(16, 19).Value = 0.182879455323625
$ws.Cells.Item(16, 20).Value = 0.182879455323625

$ws.Cells.Item(17, 7).Value = 101.145495
$ws.Cells.Item(17, 8).Value = 303.436485
$ws.Cells.Item(17, 9).Value = 0.9465500804006033
$ws.Cells.Item(17, 10).Value = 0.9465500804006032
$ws.Cells.Item(17, 13).Value = 11.78107633333333
$ws.Cells.Item(17, 14).Value = 35.343229
$ws.Cells.Item(17, 15).Value = 0.03307529931894448
$ws.Cells.Item(17, 16).Value = 0.03307529931894448
$ws.Cells.Item(17, 17).Value = 1191.602797367785
$ws.Cells.Item(17, 18).Value = 10724.42517631007
$ws.Cells.Item(17, 19).Value = 0.03130742722962091
$ws.Cells.Item(17, 20).Value = 0.03130742722962091

$ws.Cells.Item(18, 7).Value = 101.145495
$ws.Cells.Item(18, 8).Value = 303.436485
$ws.Cells.Item(18, 9).Value = 0.9465500804006033
$ws.Cells.Item(18, 10).Value = 0.9465500804006032
$ws.Cells.Item(18, 13).Value = 139.820737
$ws.Cells.Item(18, 14).Value = 419.462211
$ws.Cells.Item(18, 15).Value = 0.3925458588351179
$ws.Cells.Item(18, 16).Value = 0.3925458588351179
$ws.Cells.Item(18, 17).Value = 14142.23765512981
$ws.Cells.Item(18, 18).Value = 127280.1388961683
$ws.Cells.Item(18, 19).Value = 0.3715643142413047
$ws.Cells.Item(18, 20).Value = 0.3715643142413046

$ws.Cells.Item(19, 7).Value = 101.145495
$ws.Cells.Item(19, 8).Value = 303.436485
$ws.Cells.Item(19, 9).Value = 0.9465500804006033
$ws.Cells.Item(19, 10).Value = 0.9465500804006032
$ws.Cells.Item(19, 13).Value = 40.90925733333334
$ws.Cells.Item(19, 14).Value = 122.727772
$ws.Cells.Item(19, 15).Value = 0.1148524882558742
$ws.Cells.Item(19, 16).Value = 0.1148524882558742
$ws.Cells.Item(19, 17).Value = 4137.787083062381
$ws.Cells.Item(19, 18).Value = 37240.08374756142
$ws.Cells.Item(19, 19).Value = 0.1087136319928071
$ws.Cells.Item(19, 20).Value = 0.1087136319928071

$ws.Cells.Item(20, 7).Value = 2.152479666666667
$ws.Cells.Item(20, 8).Value = 6.457439
$ws.Cells.Item(20, 9).Value = 0.02014355460462176
$ws.Cells.Item(20, 10).Value = 0.02014355460462176
$ws.Cells.Item(20, 13).Value = 42.04602466666667
$ws.Cells.Item(20, 14).Value = 126.138074
$ws.Cells.Item(20, 15).Value = 0.1180439555498783
$ws.Cells.Item(20, 16).Value = 0.1180439555498783
$ws.Cells.Item(20, 17).Value = 90.50321315916511
$ws.Cells.Item(20, 18).Value = 814.528918432486
$ws.Cells.Item(20, 19).Value = 0.002377824864364517
$ws.Cells.Item(20, 20).Value = 0.002377824864364516

$ws.Cells.Item(21, 7).Value = 2.152479666666667
$ws.Cells.Item(21, 8).Value = 6.457439
$ws.Cells.Item(21, 9).Value = 0.02014355460462176
$ws.Cells.Item(21, 10).Value = 0.02014355460462176
$ws.Cells.Item(21, 15).Value = 0.1482760805823429
$ws.Cells.Item(21, 16).Value = 0.1482760805823429
$ws.Cells.Item(21, 17).Value = 113.6819048873627
$ws.Cells.Item(21, 18).Value = 1023.137143986264
$ws.Cells.Item(21, 19).Value = 0.002986807325769721
$ws.Cells.Item(21, 20).Value = 0.00298680732576972

$ws.Cells.Item(22, 7).Value = 2.152479666666667
$ws.Cells.Item(22, 8).Value = 6.457439
$ws.Cells.Item(22, 9).Value = 0.02014355460462176
$ws.Cells.Item(22, 10).Value = 0.02014355460462176
$ws.Cells.Item(22, 13).Value = 68.81807333333334
$ws.Cells.Item(22, 14).Value = 206.45422
$ws.Cells.Item(22, 15).Value = 0.1932063174578422
$ws.Cells.Item(22, 16).Value = 0.1932063174578422
$ws.Cells.Item(22, 17).Value = 148.1295035491756
$ws.Cells.Item(22, 18).Value = 1333.16553194258
$ws.Cells.Item(22, 19).Value = 0.003891862005669931
$ws.Cells.Item(22, 20).Value = 0.00389186200566993

$ws.Cells.Item(23, 7).Value = 2.152479666666667
$ws.Cells.Item(23, 8).Value = 6.457439
$ws.Cells.Item(23, 9).Value = 0.02014355460462176
$ws.Cells.Item(23, 10).Value = 0.02014355460462176
$ws.Cells.Item(23, 13).Value = 11.78107633333333
$ws.Cells.Item(23, 14).Value = 35.343229
$ws.Cells.Item(23, 15).Value = 0.03307529931894448
$ws.Cells.Item(23, 16).Value = 0.03307529931894448
$ws.Cells.Item(23, 17).Value = 25.35852725894789
$ws.Cells.Item(23, 18).Value = 228.226745330531
$ws.Cells.Item(23, 19).Value = 0.000666254097895367
$ws.Cells.Item(23, 20).Value = 0.0006662540978953669

$ws.Cells.Item(24, 7).Value = 2.152479666666667
$ws.Cells.Item(24, 8).Value = 6.457439
$ws.Cells.Item(24, 9).Value = 0.02014355460462176
$ws.Cells.Item(24, 10).Value = 0.02014355460462176
$ws.Cells.Item(24, 13).Value = 139.820737
$ws.Cells.Item(24, 14).Value = 419.462211
$ws.Cells.Item(24, 15).Value = 0.3925458588351179
$ws.Cells.Item(24, 16).Value = 0.3925458588351179
$ws.Cells.Item(24, 17).Value = 300.9612933708477
$ws.Cells.Item(24, 18).Value = 2708.651640337629
$ws.Cells.Item(24, 19).Value = 0.007907268942263342
$ws.Cells.Item(24, 20).Value = 0.007907268942263341

$ws.Cells.Item(25, 7).Value = 2.152479666666667
$ws.Cells.Item(25, 8).Value = 6.457439
$ws.Cells.Item(25, 9).Value = 0.02014355460462176
$ws.Cells.Item(25, 10).Value = 0.02014355460462176
$ws.Cells.Item(25, 13).Value = 40.90925733333334
$ws.Cells.Item(25, 14).Value = 122.727772
$ws.Cells.Item(25, 15).Value = 0.1148524882558742
$ws.Cells.Item(25, 16).Value = 0.1148524882558742
$ws.Cells.Item(25, 17).Value = 88.05634458843423
$ws.Cells.Item(25, 18).Value = 792.507101295908
$ws.Cells.Item(25, 19).Value = 0.002313537368658882
$ws.Cells.Item(25, 20).Value = 0.002313537368658881
